# Update numeric/classification data produced by a re-run of the NMR peak
# picking & multiplet analysis (H1_1D, COSY, HSQC sheets).

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# H1_1D sheet
# -----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("H1_1D")

# Row 2
$ws.Range("B2").Value = 7.737272812743949
$ws.Range("C2").Value = 0.9831644498137422
$ws.Range("E2").Value = "ddt"
$ws.Range("F2").Value = "0.648, 1.29, 7.67"

# Row 3
$ws.Range("B3").Value = 7.566881170536674
$ws.Range("C3").Value = 1.032298456963578
$ws.Range("E3").Value = "ddd"
$ws.Range("F3").Value = "1.27, 7.23, 7.62"

# Row 4
$ws.Range("B4").Value = 7.449063184395219
$ws.Range("C4").Value = 1.02813288006641
$ws.Range("F4").Value = "0.958, 7.68"

# Row 5
$ws.Range("B5").Value = 7.348598146833025
$ws.Range("C5").Value = 0.9976276316100109
$ws.Range("E5").Value = "ddq"
$ws.Range("F5").Value = "0.853, 7.18, 7.6"

# Row 6
$ws.Range("B6").Value = 3.308244907884105
$ws.Range("C6").Value = 1.018912934186621
$ws.Range("F6").Value = "0.793, 7.87, 17.1"

# Row 7
$ws.Range("B7").Value = 2.815544050454685
$ws.Range("C7").Value = 0.9721759318922879
$ws.Range("E7").Value = "m"
$ws.Range("F7").Value = ""

# Row 8
$ws.Range("B8").Value = 2.607303542092875
$ws.Range("C8").Value = 1.063308934024798
$ws.Range("E8").Value = "m"
$ws.Range("F8").Value = ""

# Row 9
$ws.Range("B9").Value = 1.965073516532478
$ws.Range("C9").Value = 1.092167043574244
$ws.Range("F9").Value = "4.58, 7.51, 13.7"

# Row 10
$ws.Range("B10").Value = 1.531828172669145
$ws.Range("C10").Value = 0.9708400807358168
$ws.Range("E10").Value = "ddq"
$ws.Range("F10").Value = "7.33, 9.11, 13.7"

# Row 11
$ws.Range("B11").Value = 1.001636089631483
$ws.Range("C11").Value = 3
# "7.41" looks numeric, so force text formatting before assignment, then
# drop back to the default style (this cell carried no style before).
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "7.41"
$ws.Range("F11").Style = "Normal"

# -----------------------------------------------------------------------
# COSY sheet
# -----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("COSY")

# Row 5
$ws.Range("B5").Value = 3.308415120635955
$ws.Range("C5").Value = 2.605142644437138
$ws.Range("D5").Value = 0.0933314710855484

# -----------------------------------------------------------------------
# HSQC sheet
# -----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("HSQC")

# Row 2
$ws.Range("B2").Value = 3.309042826572111
$ws.Range("C2").Value = 32.31545964835274
$ws.Range("D2").Value = -0.2761878371238708

# Row 3
$ws.Range("B3").Value = 2.809725261777252
$ws.Range("C3").Value = 32.3204001435485
$ws.Range("D3").Value = -0.005617320537567139

# Row 4
$ws.Range("B4").Value = 1.001847641444521
$ws.Range("C4").Value = 11.58371917677896
$ws.Range("D4").Value = 1.686318397521973

# Remove (former) rows 5-11, leaving only the 3 data rows above.
$ws.Range("A5:E11").EntireRow.Delete()
